# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# (commit: "Updated cryptos list on Wed Feb 14 03:55:03 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "49.486.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.636.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "325.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.051.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.644.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "49.460.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0947"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "268.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("E26").Value = "  -2.04%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0808"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.08%  "
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "128.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.68%  "
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  +6.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.062.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.15%  "
$ws.Range("E48").Value = "  -5.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.69%  "
